$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"), pushing the
# existing N:P ("Late", "heading", "Outstanding") columns to O:Q.
# Excel carries the formatting/width of the column to the left (M,
# "In Advance") onto the freshly inserted column, so mirror that here.
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Restore the selection Excel left the sheet in after the edit.
$ws.Range("R8").Select() | Out-Null
